# Update the EQTL3 (column D) time series on Sheet1 with refreshed values,
# and leave the selection where the author ended up (D38) after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    4  = 2346.2961342640269
    5  = 1008.0047307031887
    6  = -383.33124786745475
    7  = 2230.8222624931559
    8  = 5320.9952403656134
    9  = 5838.2169035851375
    10 = 7470.658368937854
    11 = 5493.0997048113531
    12 = 11308.504943659194
    13 = 11680.517848620657
    14 = 11587.516251420138
    15 = 10348.675117827352
    16 = 9768.6833032279974
    17 = 10155.942786401696
    18 = 10330.306705276818
    19 = 9547.3412998822914
    20 = 5116.222974985365
    21 = 5225.6606237213318
    22 = 7421.6683063923683
    23 = 8047.2674890790931
    24 = 17133.14369506764
    25 = 6640.3369228731744
    26 = -36256.629550042744
    27 = 58720.223253545359
    28 = 1953.156109231425
    29 = 2076.2237893806973
    30 = 5663.0097529148252
    31 = 2025.9664180611614
    32 = 1919.8539728117053
    33 = 1809.1392026497824
    34 = -438.17127495682212
    35 = 1515.8477492685827
    36 = 1521.1057646670865
    37 = 1525.4028087225845
    38 = 17770.841584153739
}

foreach ($row in $values.Keys) {
    $ws.Range("D$row").Value = $values[$row]
}

# Restore Sheet1 as the active sheet/selection, matching the cell the
# author was last on (D38) when the workbook was saved.
$ws.Activate()
$ws.Range("D38").Select()
